$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.94"
$ws.Range("E2").Value = "'-4.52%"
$ws.Range("D3").Value = "'40.02"
$ws.Range("E3").Value = "'-6.95%"
$ws.Range("D4").Value = "'5.147"
$ws.Range("E4").Value = "'-1.10%"
$ws.Range("D5").Value = "'0.07699"
$ws.Range("E5").Value = "'-6.08%"
$ws.Range("D6").Value = "'4.264"
$ws.Range("E6").Value = "'-1.27%"
$ws.Range("D7").Value = "'1.627"
$ws.Range("E7").Value = "'-11.90%"
$ws.Range("D8").Value = "'0.8784"
$ws.Range("E8").Value = "'-6.47%"
$ws.Range("D9").Value = "'0.09962"
$ws.Range("E9").Value = "'-10.60%"
$ws.Range("E10").Value = "'-6.89%"
$ws.Range("D11").Value = "'0.08936"
$ws.Range("E11").Value = "'-4.43%"
$ws.Range("D12").Value = "'0.04407"
$ws.Range("E12").Value = "'-4.46%"
$ws.Range("E13").Value = "'-0.43%"
$ws.Range("D14").Value = "'0.001250"
$ws.Range("E14").Value = "'-2.52%"
$ws.Range("D15").Value = "'0.005954"
$ws.Range("E15").Value = "'-1.57%"
$ws.Range("D16").Value = "'3.355"
$ws.Range("E16").Value = "'-0.06%"
$ws.Range("D18").Value = "'0.3321"
$ws.Range("E18").Value = "'-1.34%"
$ws.Range("D19").Value = "'6.981"
$ws.Range("E19").Value = "'-5.98%"
$ws.Range("E20").Value = "'-3.31%"
$ws.Range("D21").Value = "'0.3134"
$ws.Range("D22").Value = "'0.04145"
$ws.Range("E22").Value = "'-0.23%"
$ws.Range("D23").Value = "'0.001198"
$ws.Range("E23").Value = "'-4.00%"
$ws.Range("D24").Value = "'0.004069"
$ws.Range("E24").Value = "'-5.63%"
$ws.Range("D25").Value = "'0.0001219"
$ws.Range("E25").Value = "'10.94%"
$ws.Range("E26").Value = "'0.19%"
$ws.Range("E38").Value = "'-14.06%"
$ws.Range("D39").Value = "'0.05153"
$ws.Range("E39").Value = "'-7.06%"
$ws.Range("D40").Value = "'0.007939"
$ws.Range("E40").Value = "'-0.80%"
$ws.Range("D41").Value = "'0.1323"
$ws.Range("E41").Value = "'-5.22%"
$ws.Range("D42").Value = "'0.006385"
$ws.Range("E42").Value = "'-2.42%"
$ws.Range("D43").Value = "'0.001941"
$ws.Range("E43").Value = "'-8.38%"
$ws.Range("D44").Value = "'0.008597"
$ws.Range("E44").Value = "'3.71%"
$ws.Range("D45").Value = "'0.3054"
$ws.Range("E45").Value = "'-4.76%"
$ws.Range("D46").Value = "'0.00006500"
$ws.Range("E46").Value = "'-6.32%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'0.04%"
$ws.Range("D48").Value = "'0.006995"
$ws.Range("E48").Value = "'98.38%"
$ws.Range("D49").Value = "'0.003393"
$ws.Range("E49").Value = "'-1.99%"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'0.04%"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'0.04%"
